# Necron_units.xlsx edit script
# Adds a new "Options" column (E) to every sheet, fills in the new wargear
# "Options" alternatives that have been worked out so far, fixes a
# Lychguard/Triarch Praetorians wargear swap on the Elites sheet, and backfills
# the missing "Wargear" header on the Heavy Support sheet.
# Finishes with "Elites" as the active sheet/tab (matching the workbook's
# new activeTab), each sheet keeping its own last-used selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Named Characters
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Named Characters")
$ws.Activate()
$ws.Range("E1").Value = "Options"
$ws.Range("E1").Select()

# ---------------------------------------------------------------------------
# HQ
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HQ")
$ws.Activate()
$ws.Range("E1").Value = "Options"
$ws.Range("E2").Value = "Tesla Cannon, Staff of light/Warscythe/Hyperphase sword/Voidblade/, Phylactery, Resurrection orb"
$ws.Range("E3").Value = "Phylactery, Canoptek Cloak/Chronometron"
$ws.Range("E4").Value = "Staff of light/Warscythe/Hyperphase sword/Voidblade/, Phylactery, Resurrection orb"
$ws.Range("E5").Value = "Staff of light/Warscythe/Hyperphase sword/Voidblade/, Phylactery, Resurrection orb"
$ws.Range("E6").Value = "Staff of light/Warscythe/Hyperphase sword/Voidblade/, Phylactery, Resurrection orb"
$ws.Columns.Item(4).ColumnWidth = 24.451822916666668
$ws.Range("E10").Select()

# ---------------------------------------------------------------------------
# Troops
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Troops")
$ws.Activate()
$ws.Range("E1").Value = "Options"
$ws.Columns.Item(4).ColumnWidth = 11.877604166666666
$ws.Range("E1").Select()

# ---------------------------------------------------------------------------
# Elites
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elites")
$ws.Activate()
$ws.Range("E1").Value = "Options"
# Lychguard / Triarch Praetorians wargear swap
$ws.Range("D4").Value = "Warscythe"
$ws.Range("D5").Value = "Rod of covenant"
$ws.Range("E4").Value = "Warscythe/Hyperphase sword+Dispersion shield"
$ws.Range("E5").Value = "Rod of the covennant/Particle caster+Voidblade"
$ws.Range("E6").Value = "Heat Ray/2*Heavy gauss cannon/Particle Shredder"
$ws.Columns.Item(4).ColumnWidth = 25.166666666666668
$ws.Range("D5").Select()

# ---------------------------------------------------------------------------
# Fast Attack
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Fast Attack")
$ws.Activate()
$ws.Range("E1").Value = "Options"
$ws.Range("E3").Value = "Whip coils, Transdimensional beamer/Particle caster"
$ws.Range("E4").Value = "Heavy gauss cannon-3"
$ws.Range("E5").Value = "Shieldvanes, Nebuloscope/Shadowloom"
$ws.Columns.Item(4).ColumnWidth = 16.451822916666668
$ws.Range("E7").Select()

# ---------------------------------------------------------------------------
# Heavy Support  (this sheet was missing the "Wargear" header altogether)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Heavy Support")
$ws.Activate()
$ws.Range("D1").Value = "Wargear"
$ws.Range("E1").Value = "Options"
$ws.Range("E2").Value = "Gauss cannon/Tesla cannon"
$ws.Range("E3").Value = "Gloom prism, Fabricator claw array"
$ws.Columns.Item(3).ColumnWidth = 15.451822916666666
$ws.Columns.Item(4).ColumnWidth = 35.736979166666664
$ws.Range("E7").Select()

# ---------------------------------------------------------------------------
# Dedicated Transports
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dedicated Transports")
$ws.Activate()
$ws.Range("E1").Value = "Options"
$ws.Columns.Item(4).ColumnWidth = 18.022135416666668
$ws.Range("E2").Select()

# ---------------------------------------------------------------------------
# Flyers
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Flyers")
$ws.Activate()
$ws.Range("E1").Value = "Options"
$ws.Range("E1").Select()

# ---------------------------------------------------------------------------
# Final active sheet/tab: Elites (workbook activeTab index 3)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elites")
$ws.Activate()
$ws.Range("D5").Select()
